$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume columns so values like "51.50" or
# "0.08230" keep their trailing zeros instead of being parsed as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.957.70'
$ws.Range("E2").Value = '  -2.91%  '

$ws.Range("D3").Value = '1.886.99'
$ws.Range("E3").Value = '  -3.86%  '

$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -1.02%  '

$ws.Range("D5").Value = '325.73'
$ws.Range("E5").Value = '  +0.82%  '

$ws.Range("E6").Value = '  -0.70%  '

$ws.Range("D7").Value = '0.4574'
$ws.Range("E7").Value = '  -4.14%  '

$ws.Range("D8").Value = '0.3927'
$ws.Range("E8").Value = '  -2.76%  '

$ws.Range("D9").Value = '51.50'
$ws.Range("E9").Value = '  -4.48%  '

$ws.Range("D10").Value = '0.08230'
$ws.Range("E10").Value = '  -2.75%  '

$ws.Range("D11").Value = '1.034'
$ws.Range("E11").Value = '  -2.50%  '

$ws.Range("E12").Value = '  -3.81%  '

$ws.Range("D13").Value = '1.907.09'
$ws.Range("E13").Value = '  -3.43%  '

$ws.Range("D14").Value = '7.309'
$ws.Range("E14").Value = '  -4.45%  '

$ws.Range("D15").Value = '5.973'
$ws.Range("E15").Value = '  -4.33%  '

$ws.Range("E16").Value = '  -0.85%  '

$ws.Range("D17").Value = '89.04'
$ws.Range("E17").Value = '  -0.87%  '

$ws.Range("D18").Value = '0.00001056'
$ws.Range("E18").Value = '  -1.16%  '

$ws.Range("E19").Value = '  -0.46%  '

$ws.Range("D20").Value = '17.56'
$ws.Range("E20").Value = '  -5.96%  '

$ws.Range("E21").Value = '  -0.68%  '

$ws.Range("D22").Value = '5.640'
$ws.Range("E22").Value = '  -2.53%  '

$ws.Range("D23").Value = '27.979.46'
$ws.Range("E23").Value = '  -2.90%  '

$ws.Range("E24").Value = '  -3.88%  '

$ws.Range("D25").Value = '2.305'
$ws.Range("E25").Value = '  +0.92%  '

$ws.Range("D26").Value = '2.156.72'
$ws.Range("E26").Value = '  -2.07%  '

$ws.Range("D27").Value = '154.11'
$ws.Range("E27").Value = '  -0.21%  '

$ws.Range("D28").Value = '19.85'
$ws.Range("E28").Value = '  -1.80%  '

$ws.Range("D29").Value = '2.102'
$ws.Range("E29").Value = '  -2.61%  '

$ws.Range("E30").Value = '  -4.85%  '

$ws.Range("D31").Value = '123.99'
$ws.Range("E31").Value = '  -0.16%  '

$ws.Range("D32").Value = '0.09522'
$ws.Range("E32").Value = '  -0.91%  '

$ws.Range("D33").Value = '0.9559'
$ws.Range("E33").Value = '  -4.88%  '

$ws.Range("E34").Value = '  +0.21%  '

$ws.Range("D35").Value = '3.631'
$ws.Range("E35").Value = '  -1.36%  '

$ws.Range("D36").Value = '5.460'
$ws.Range("E36").Value = '  -3.80%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '1.253'
$ws.Range("E37").Value = '  -1.12%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02275'
$ws.Range("E38").Value = '  -3.45%  '

$ws.Range("D39").Value = '8.642'
$ws.Range("E39").Value = '  -1.24%  '

$ws.Range("D40").Value = '0.06097'
$ws.Range("E40").Value = '  -1.78%  '

$ws.Range("D41").Value = '0.6093'
$ws.Range("E41").Value = '  -2.34%  '

$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  -0.62%  '

$ws.Range("D43").Value = '10.71'
$ws.Range("E43").Value = '  -3.65%  '

$ws.Range("D44").Value = '0.1883'
$ws.Range("E44").Value = '  -1.80%  '

$ws.Range("D45").Value = '1.306'
$ws.Range("E45").Value = '  -2.87%  '

$ws.Range("D46").Value = '0.5803'
$ws.Range("E46").Value = '  -2.61%  '

$ws.Range("D47").Value = '12.69'
$ws.Range("E47").Value = '  -1.99%  '

$ws.Range("D48").Value = '1.987'
$ws.Range("E48").Value = '  -4.50%  '

$ws.Range("D49").Value = '3.422'
$ws.Range("E49").Value = '  +0.29%  '

$ws.Range("D50").Value = '0.06880'
$ws.Range("E50").Value = '  +0.58%  '

$ws.Range("D51").Value = '110.11'
$ws.Range("E51").Value = '  -0.70%  '
